$wb = $excel.ActiveWorkbook

# --- Sources sheet: insert a new row into SourceTable (ADD101 / ptntid) ---
$wsSrc = $wb.Worksheets.Item("Sources")
$wsSrc.Rows(15).Insert()
$loSrc = $wsSrc.ListObjects.Item("SourceTable")
$loSrc.Resize($wsSrc.Range("A1:E17"))
$wsSrc.Range("A14:E14").Copy()
$wsSrc.Range("A15:E15").PasteSpecial(-4122)
$wsSrc.Range("A15").Value = "ADD101"
$wsSrc.Range("B15").Value = "ptntid"
$wsSrc.Range("C15").Value = "int"
$wsSrc.Range("D15").Value = "1,2"
$wsSrc.Range("E15").Value = "N"

# --- Target sheet: insert a new row into TargetTable (InsuranceName / latestInsurance) ---
# and rename the 4th column header from TargetColumnType to ColumnDataType
$wsTgt = $wb.Worksheets.Item("Target")
$wsTgt.Rows(4).Insert()
$loTgt = $wsTgt.ListObjects.Item("TargetTable")
$loTgt.Resize($wsTgt.Range("A1:D12"))
$wsTgt.Range("A3:D3").Copy()
$wsTgt.Range("A4:D4").PasteSpecial(-4122)
$wsTgt.Range("A4").Value = "PATIENT"
$wsTgt.Range("B4").Value = "InsuranceName"
$wsTgt.Range("C4").Value = "latestInsurance"
$wsTgt.Range("D4").Value = "varchar"
$wsTgt.Range("D1").Value = "ColumnDataType"

# --- Selection / active tab restore (Sources becomes the active tab) ---
$wsTgt.Range("D2").Select()
$wsSrc.Activate()
$wsSrc.Range("C1").Select()
